$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$codes = $wb.Worksheets.Item("Include #0")

# URL row
$metadata.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-medication-reconciliation-resolution"

# Name row
$metadata.Range("B4").Value = "FRMedicationReconciliationResolution"

# Title row
$metadata.Range("B5").Value = "value set Interop'Santé - Résolution d'une divergence sur une ligne de traitement d'une FCT"

# Date row
$metadata.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction row
$metadata.Range("B11").Value = "FRANCE"

# System URI row on the codes sheet
$codes.Range("B4").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-reconciliation-resolution"
